$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.798.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Value = "'2.047.79"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +0.88%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("E4").Value = "'  -0.02%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'227.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  +0.18%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.32%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'59.89"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.70%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "'  +0.01%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'  -1.70%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'  +3.12%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("E11").Value = "'  +0.07%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'2.351.57"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  +0.82%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'14.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.36%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'21.45"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +1.53%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'5.51"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +6.37%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("E16").Value = "'  +0.20%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'2.052.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  +1.09%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'37.809.63"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +0.11%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'69.56"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.54%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'5.92"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  -1.80%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'0.0₃0830"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  +0.67%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'222.50"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -1.04%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("E23").Value = "'  +0.45%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'  -2.29%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  +2.87%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'168.91"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +2.23%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'9.33"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  +0.51%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'  -0.82%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "'  -0.78%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = "'  +0.58%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("E31").Value = "'  -0.39%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'2.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +8.40%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'4.39"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -1.01%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("E34").Value = "'  +0.82%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'0.0603"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +0.11%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'6.49"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +2.03%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("E37").Value = "'  +4.20%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  +7.70%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("E39").Value = "'  -0.06%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'18.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +9.50%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'1.525.87"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.01%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'97.73"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.19%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'  -1.20%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("E44").Value = "'  +0.00%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'  +1.24%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -2.90%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "'  -0.27%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("E48").Value = "'  +0.36%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'7.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  +0.22%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -0.50%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'2.239.22"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +0.85%  "
$ws.Range("E51").Style = "Normal"
